$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Limiti" - header text + updated monthly limits / remaining amounts
# ---------------------------------------------------------------------------
$wsLimiti = $wb.Worksheets.Item("Limiti")

$wsLimiti.Range("A1").Value = "Kategorija"

$wsLimiti.Range("B2").Value = 50
$wsLimiti.Range("C2").Value = 35.68

$wsLimiti.Range("B3").Value = 12
$wsLimiti.Range("C3").Value = 8.5

$wsLimiti.Range("B4").Value = 30
$wsLimiti.Range("C4").Value = 16.79

$wsLimiti.Range("B5").Value = 12
$wsLimiti.Range("C5").Value = 0

$wsLimiti.Range("B6").Value = 30
$wsLimiti.Range("C6").Value = 16.79

$wsLimiti.Range("B7").Value = 12
$wsLimiti.Range("C7").Value = 12

$wsLimiti.Range("B8").Value = 30
$wsLimiti.Range("C8").Value = 30

$wsLimiti.Range("B9").Value = 60
$wsLimiti.Range("C9").Value = 0

$wsLimiti.Range("B10").Value = 1
$wsLimiti.Range("C10").Value = 1

$wsLimiti.Range("B11").Value = 15
$wsLimiti.Range("C11").Value = 15

$wsLimiti.Range("B12").Value = 21
$wsLimiti.Range("C12").Value = 21

$wsLimiti.Range("B13").Value = 23
$wsLimiti.Range("C13").Value = 0

# ---------------------------------------------------------------------------
# Sheet "Izdevumi" - new expense log entries (rows 2-9)
# ---------------------------------------------------------------------------
$wsIzdevumi = $wb.Worksheets.Item("Izdevumi")

$wsIzdevumi.Range("A2").Value = "2025-05-20 17:52:39"
$wsIzdevumi.Range("B2").Value = "pārtika"
$wsIzdevumi.Range("C2").Value = 14.32

$wsIzdevumi.Range("A3").Value = "2025-05-20 17:52:48"
$wsIzdevumi.Range("B3").Value = "mājas izdevumi (komunālie + īre / nekustamā īpašuma nodoklis)"
$wsIzdevumi.Range("C3").Value = 13.21

$wsIzdevumi.Range("A4").Value = "2025-05-20 17:52:56"
$wsIzdevumi.Range("B4").Value = "hobiji"
$wsIzdevumi.Range("C4").Value = 12.42

$wsIzdevumi.Range("A5").Value = "2025-05-20 17:53:10"
$wsIzdevumi.Range("B5").Value = "ēšana ārpus mājas (restorāni/fast food/kafejnīcas)"
$wsIzdevumi.Range("C5").Value = 3.5

$wsIzdevumi.Range("A6").Value = "2025-05-20 17:53:16"
$wsIzdevumi.Range("B6").Value = "dāvanas"
$wsIzdevumi.Range("C6").Value = 11.34

$wsIzdevumi.Range("A7").Value = "2025-05-20 17:53:32"
$wsIzdevumi.Range("B7").Value = "dāvanas"
$wsIzdevumi.Range("C7").Value = 13.21

$wsIzdevumi.Range("A8").Value = "2025-05-20 17:53:42"
$wsIzdevumi.Range("B8").Value = "medicīniskie izdevumi"
$wsIzdevumi.Range("C8").Value = 120

$wsIzdevumi.Range("A9").Value = "2025-05-20 17:53:53"
$wsIzdevumi.Range("B9").Value = "mājdzīvnieki"
$wsIzdevumi.Range("C9").Value = 13.21
